# Injector assembly workbook - add a "Calcs" worksheet with injector
# sizing labels (discharge coefficient, areas, densities) after Sheet1.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after Sheet1 and rename it.
$calcs = $wb.Worksheets.Add($null, $sheet1)
$calcs.Name = "Calcs"

# Title in A1, bigger font, matching row height.
$calcs.Range("A1").Value = "Injector"
$calcs.Range("A1").Font.Name = "Calibri (Body)"
$calcs.Range("A1").Font.Size = 20
$calcs.Rows.Item(1).RowHeight = 26

# Labels down column A (shared-string table is populated in authoring order).
$calcs.Range("A3").Value = "Discharge Coefficient (Cd) "
$calcs.Range("A13").Value = "Oxidizer cross sectional area (Ao)"
$calcs.Range("A15").Value = "Fuel cross sectional area (Af)"
$calcs.Range("A5").Value = "Oxidizer Density (Do)"
$calcs.Range("A7").Value = "Fuel Density (Df)"

# Center-aligned, merged input cells next to each label.
$calcs.Range("C3:E3").HorizontalAlignment = -4108
$calcs.Range("C5:E5").HorizontalAlignment = -4108
$calcs.Range("C7:E7").HorizontalAlignment = -4108
$calcs.Range("C3:E3").Merge()
$calcs.Range("C5:E5").Merge()
$calcs.Range("C7:E7").Merge()

# Widen column A to fit the longest label.
$calcs.Columns.Item(1).ColumnWidth = 28.166666666666668

# Leave the cursor where the author left it and make Calcs the active tab.
[void]$calcs.Range("A11").Select()
[void]$calcs.Activate()
